# CPT-280 style: customize header row's stylings
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

# --- Text updates: append "Locale: " prefix to the zh/en language labels ---
$ws.Range("B1").Value = "Locale: zh"
$ws.Range("C1").Value = "Locale: en"

# --- Make the header row taller ---
$ws.Rows.Item(1).RowHeight = 40

# --- Column A keeps its width; widen columns B and C ---
$ws.Columns.Item(2).ColumnWidth = 149.17
$ws.Columns.Item(3).ColumnWidth = 149.17

# --- Center (horizontally & vertically) all header cells ---
$headerRange = $ws.Range("A1:C1")
$headerRange.HorizontalAlignment = $xlCenter
$headerRange.VerticalAlignment = $xlCenter

# --- Emphasize the language-label cells: bold + brand blue font color ---
$labelRange = $ws.Range("B1:C1")
$labelRange.Font.Color = 15351613
$labelRange.Font.Bold = $true

# --- Freeze both the header row and the first column, matching the new B2 split ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B2").Select()
$win.FreezePanes = $true
